$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for columns B:G, rows 2-17 (HIGH, LOW, CLOSE, LTP, VOL, 9:25 CLOSE)
$values = @(
    @(1434.35, 1412.1,  3136,     1415.15, 36,  1423.8),
    @(1505.25, 1486.4,  1499.45,  1499.25, 21,  1493),
    @(51818.25,51219.1, 51391.05, 51386.4, 16,  51630.15),
    @(105.26,  104.06,  104.35,   104.45,  389, 104.62),
    @(864.25,  846,     851.1,    851.35,  51,  858.8),
    @(752.05,  740.05,  749.8,    749.35,  117, 741.15),
    @(1242.75, 1227.3,  1228.7,   1229,    184, 1240.35),
    @(1017,    999.45,  1009.2,   1008.95, 30,  1000.3),
    @(25150,   25001.1, 25048,    25049.55,32,  25089),
    @(2782.15, 2750.9,  2758.55,  2755.35, 78,  2767.35),
    @(806.7,   795.8,   802.4,    803.35,  200, 800.75),
    @(1122.85, 1105.15, 1118,     1118.1,  16,  1117.7),
    @(941.45,  931.75,  934.45,   933.4,   58,  934.7),
    @(164.7,   160.88,  161.56,   161.36,  713, 161.46),
    @(4252.7,  4122.9,  4157,     4159.55, 97,  4208.55),
    @(3504.9,  3468.6,  3484.65,  3489.6,  15,  3491.55)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
